## Weekly update: prepend a new price record for "Alcachofa" (Macroferia
## Regional de Talca) as row 56, pushing the existing historical rows
## (56-111) down by one (now 57-112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row just above the current row 56; this shifts all
# rows 56..111 down to 57..112 and grows the used range to A1:R112.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A56").Value = 5
$ws.Range("B56").Value = "Macroferia Regional de Talca"
$ws.Range("C56").Value = "Maule"
$ws.Range("D56").Value = 45033
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = 100112013
$ws.Range("G56").Value = "Alcachofa"
$ws.Range("H56").Value = "Argentina(o)"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 200
$ws.Range("K56").Value = 17000
$ws.Range("L56").Value = 17000
$ws.Range("M56").Value = 17000
$ws.Range("N56").Value = "$/caja 40 unidades"
$ws.Range("O56").Value = "Provincia del Elquí"
$ws.Range("P56").Value = 425
$ws.Range("Q56").Value = 40
$ws.Range("R56").Value = "Hortaliza"
